# Update the RandomForestRegressor prediction row (row 2) with new metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1026899336953756
$ws.Range("C2").Value = 0.9993883053647487
$ws.Range("D2").Value = 0.2199866449939377
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=7))])"
$ws.Range("G2").Value = 0.12469127785007
$ws.Range("H2").Value = 0.99
